$d = $word.ActiveDocument

# 1. Update the letter date from September 19, 2025 to September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the sender/return address line (the one right under "Kevin Nutter",
#    NOT the later "PROPERTY ADDRESS:" one) into two paragraphs:
#    "2191 Rancho Mccormick Blvd., Santa Clara CA 95050"
# -> "2191 Rancho Mccormick Blvd."
#    "Santa Clara, CA 95050"
$targetText = "2191 Rancho Mccormick Blvd., Santa Clara CA 95050"
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq $targetText) {
        # Replace the whole paragraph range (which includes the trailing
        # paragraph mark) with the street portion, a paragraph break, and the
        # city/state/zip portion in one shot. This preserves the paragraph
        # mark/run formatting (Arial, 11pt) on both resulting paragraphs and
        # avoids leaving behind a stray empty run.
        $p.Range.Text = "2191 Rancho Mccormick Blvd." + [char]13 + "Santa Clara, CA 95050"
        break
    }
}

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "... Board of Directors" line.
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Rancho Palma Grande Homeowners Association Board of Directors") {
        $nextPara = $p.Next()
        $nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)
        if ($nextText -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
